$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CapitalDistribution")

# Rename the fund from "Agri Fund" to "SAAS Fund" across all data rows
$ws.Range("A2").Value = "SAAS Fund"
$ws.Range("A3").Value = "SAAS Fund"
$ws.Range("A4").Value = "SAAS Fund"

# Update the active selection to match the edited state
$ws.Range("A4").Select()
